$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.996.78"
$ws.Range("E2").Value = "  +2.25%  "

# Row 3
$ws.Range("D3").Value = "2.231.54"
$ws.Range("E3").Value = "  +1.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.42"
$ws.Range("E5").Value = "  -1.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.15"
$ws.Range("E6").Value = "  +5.59%  "

# Row 7
$ws.Range("E7").Value = "  +1.56%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  +1.48%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.09"
$ws.Range("E10").Value = "  +7.62%  "

# Row 11
$ws.Range("E11").Value = "  +2.44%  "

# Row 12
$ws.Range("E12").Value = "  -0.45%  "

# Row 13
$ws.Range("E13").Value = "  +1.34%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.39"
$ws.Range("E14").Value = "  +2.57%  "

# Row 15
$ws.Range("D15").Value = "2.582.77"
$ws.Range("E15").Value = "  +1.78%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.05"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("D17").Value = "2.248.98"
$ws.Range("E17").Value = "  +2.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.730"
$ws.Range("E18").Value = "  +2.71%  "

# Row 19
$ws.Range("D19").Value = "39.920.33"
$ws.Range("E19").Value = "  +2.35%  "

# Row 20
$ws.Range("E20").Value = "  +2.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("E21").Value = "  +8.23%  "

# Row 22
$ws.Range("E22").Value = "  +2.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.45"
$ws.Range("E23").Value = "  +1.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.51"
$ws.Range("E24").Value = "  +3.36%  "

# Row 25
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("E26").Value = "  +3.20%  "

# Row 27
$ws.Range("E27").Value = "  +4.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.80"
$ws.Range("E28").Value = "  +1.76%  "

# Row 29
$ws.Range("E29").Value = "  +3.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.31"
$ws.Range("E30").Value = "  +2.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.11"
$ws.Range("E31").Value = "  +4.65%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.95"
$ws.Range("E32").Value = "  +1.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  +2.16%  "

# Row 35
$ws.Range("E35").Value = "  +4.40%  "

# Row 36
$ws.Range("E36").Value = "  +2.86%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.23"
$ws.Range("E37").Value = "  +8.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("E38").Value = "  +7.75%  "

# Row 39
$ws.Range("E39").Value = "  +2.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1000"
$ws.Range("E40").Value = "  +4.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.71"
$ws.Range("E41").Value = "  +6.46%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.81"
$ws.Range("E42").Value = "  +4.70%  "

# Row 43
$ws.Range("D43").Value = "2.067.09"
$ws.Range("E43").Value = "  +8.58%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.13"
$ws.Range("E44").Value = "  +14.18%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.12"
$ws.Range("E45").Value = "  +5.39%  "

# Row 46
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("E46").Value = "  +4.57%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.90"
$ws.Range("E47").Value = "  +11.40%  "

# Row 48
$ws.Range("E48").Value = "  -0.80%  "

# Row 49
$ws.Range("D49").Value = "2.438.51"
$ws.Range("E49").Value = "  +1.13%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.67"
$ws.Range("E50").Value = "  +2.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.09"
$ws.Range("E51").Value = "  +3.49%  "
